# Se refactoriza la tarea DescargarSaldo
#
# Updates the "Datos" sheet (first/active sheet) test-data row used by the
# DescargarSaldo data-driven task:
#   - B2 (numeroDocumento): 333333304 -> 93221451
#   - D2 (tipoDocumento):   "autotest32" -> "autotest26"
# and leaves the active selection on B2 (matching the last-edited cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 93221451
$ws.Range("D2").Value = "autotest26"

$ws.Range("B2").Select() | Out-Null
